# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.266.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "'1.683.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'217.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'0.5247"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "'0.2702"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("D9").Value = "'0.06408"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "'21.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("D11").Value = "'0.07491"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").Value = "'1.704.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "'0.000008446"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("D16").Value = "'64.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "'26.313.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'4.920"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'189.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").Value = "'6.192"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").Value = "'1.008"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'144.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'7.694"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  +4.44%  "
$ws.Range("D27").Value = "'15.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").Value = "'0.06648"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.55%  "
$ws.Range("D29").Value = "'1.347"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.31%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").Value = "'3.565"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  +1.32%  "
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("D38").Value = "'6.387"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'0.01621"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "'1.104.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("D41").Value = "'0.8753"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "'1.014"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").Value = "'100.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").Value = "'1.832.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").Value = "'56.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "'8.146"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("D49").Value = "'0.05271"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").Value = "'0.4303"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").Value = "'6.021"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.26%  "
